# Update "想去人数" (want-to-go count) values for several manga-expo events,
# reflecting refreshed data at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value = 9496   # 合肥·第六届环形宇宙动漫游戏嘉年华-一周年超强巨制~
$wsExpo.Range("F7").Value = 854    # 合肥·第二届华盟动漫次元嘉年华
$wsExpo.Range("F10").Value = 1180  # 合肥·城市动漫节

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 15     # 合肥·跨越二次元ACG神级动漫世界巡回演唱会

# Sheet "全部类型" (all types combined)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 15      # 合肥·跨越二次元ACG神级动漫世界巡回演唱会
$wsAll.Range("F7").Value = 9496    # 合肥·第六届环形宇宙动漫游戏嘉年华-一周年超强巨制~
$wsAll.Range("F8").Value = 854     # 合肥·第二届华盟动漫次元嘉年华
$wsAll.Range("F11").Value = 1180   # 合肥·城市动漫节
